# Natmi following Dr Hou advice
# Recompute the Tfpi-Vldlr ligand-receptor table over the full 4x4 grid of
# sending/target clusters (ECs, FAPs, M2, sCs) -- adds the previously
# missing "M2" target-cluster rows and refreshes every numeric column with
# the revised per-cluster denominators.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tfpi"
$ws.Cells.Item(2, 3).Value = "Vldlr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 5.451731666666667
$ws.Cells.Item(2, 8).Value = 16.355195
$ws.Cells.Item(2, 9).Value = 0.1922099906071488
$ws.Cells.Item(2, 10).Value = 0.1922099906071488
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.463191
$ws.Cells.Item(2, 14).Value = 1.389573
$ws.Cells.Item(2, 15).Value = 0.0353316468093919
$ws.Cells.Item(2, 16).Value = 0.0353316468093919
$ws.Cells.Item(2, 17).Value = 2.525193042415
$ws.Cells.Item(2, 18).Value = 22.726737381735
$ws.Cells.Item(2, 19).Value = 0.006791095501368317
$ws.Cells.Item(2, 20).Value = 0.006791095501368317

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tfpi"
$ws.Cells.Item(3, 3).Value = "Vldlr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 5.451731666666667
$ws.Cells.Item(3, 8).Value = 16.355195
$ws.Cells.Item(3, 9).Value = 0.1922099906071488
$ws.Cells.Item(3, 10).Value = 0.1922099906071488
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.932929333333332
$ws.Cells.Item(3, 14).Value = 29.798788
$ws.Cells.Item(3, 15).Value = 0.7576717833204485
$ws.Cells.Item(3, 16).Value = 0.7576717833204486
$ws.Cells.Item(3, 17).Value = 54.15166538929556
$ws.Cells.Item(3, 18).Value = 487.36498850366
$ws.Cells.Item(3, 19).Value = 0.1456320863553251
$ws.Cells.Item(3, 20).Value = 0.1456320863553251

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tfpi"
$ws.Cells.Item(4, 3).Value = "Vldlr"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 5.451731666666667
$ws.Cells.Item(4, 8).Value = 16.355195
$ws.Cells.Item(4, 9).Value = 0.1922099906071488
$ws.Cells.Item(4, 10).Value = 0.1922099906071488
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1285543333333333
$ws.Cells.Item(4, 14).Value = 0.385663
$ws.Cells.Item(4, 15).Value = 0.009805968382697785
$ws.Cells.Item(4, 16).Value = 0.009805968382697785
$ws.Cells.Item(4, 17).Value = 0.7008437299205555
$ws.Cells.Item(4, 18).Value = 6.307593569285
$ws.Cells.Item(4, 19).Value = 0.00188480509073234
$ws.Cells.Item(4, 20).Value = 0.00188480509073234

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Tfpi"
$ws.Cells.Item(5, 3).Value = "Vldlr"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.451731666666667
$ws.Cells.Item(5, 8).Value = 16.355195
$ws.Cells.Item(5, 9).Value = 0.1922099906071488
$ws.Cells.Item(5, 10).Value = 0.1922099906071488
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.585130333333333
$ws.Cells.Item(5, 14).Value = 7.755391
$ws.Cells.Item(5, 15).Value = 0.1971906014874617
$ws.Cells.Item(5, 16).Value = 0.1971906014874618
$ws.Cells.Item(5, 17).Value = 14.09343690069389
$ws.Cells.Item(5, 18).Value = 126.840932106245
$ws.Cells.Item(5, 19).Value = 0.03790200365972304
$ws.Cells.Item(5, 20).Value = 0.03790200365972305

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tfpi"
$ws.Cells.Item(6, 3).Value = "Vldlr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 15.797976
$ws.Cells.Item(6, 8).Value = 47.393928
$ws.Cells.Item(6, 9).Value = 0.5569842765993244
$ws.Cells.Item(6, 10).Value = 0.5569842765993244
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.463191
$ws.Cells.Item(6, 14).Value = 1.389573
$ws.Cells.Item(6, 15).Value = 0.0353316468093919
$ws.Cells.Item(6, 16).Value = 0.0353316468093919
$ws.Cells.Item(6, 17).Value = 7.317480301415999
$ws.Cells.Item(6, 18).Value = 65.85732271274401
$ws.Cells.Item(6, 19).Value = 0.01967917173919197
$ws.Cells.Item(6, 20).Value = 0.01967917173919197

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tfpi"
$ws.Cells.Item(7, 3).Value = "Vldlr"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 15.797976
$ws.Cells.Item(7, 8).Value = 47.393928
$ws.Cells.Item(7, 9).Value = 0.5569842765993244
$ws.Cells.Item(7, 10).Value = 0.5569842765993244
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 9.932929333333332
$ws.Cells.Item(7, 14).Value = 29.798788
$ws.Cells.Item(7, 15).Value = 0.7576717833204485
$ws.Cells.Item(7, 16).Value = 0.7576717833204486
$ws.Cells.Item(7, 17).Value = 156.920179217696
$ws.Cells.Item(7, 18).Value = 1412.281612959264
$ws.Cells.Item(7, 19).Value = 0.42201127013246
$ws.Cells.Item(7, 20).Value = 0.4220112701324601

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Tfpi"
$ws.Cells.Item(8, 3).Value = "Vldlr"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 15.797976
$ws.Cells.Item(8, 8).Value = 47.393928
$ws.Cells.Item(8, 9).Value = 0.5569842765993244
$ws.Cells.Item(8, 10).Value = 0.5569842765993244
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1285543333333333
$ws.Cells.Item(8, 14).Value = 0.385663
$ws.Cells.Item(8, 15).Value = 0.009805968382697785
$ws.Cells.Item(8, 16).Value = 0.009805968382697785
$ws.Cells.Item(8, 17).Value = 2.030898272696
$ws.Cells.Item(8, 18).Value = 18.278084454264
$ws.Cells.Item(8, 19).Value = 0.005461770205992773
$ws.Cells.Item(8, 20).Value = 0.005461770205992773

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Tfpi"
$ws.Cells.Item(9, 3).Value = "Vldlr"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 15.797976
$ws.Cells.Item(9, 8).Value = 47.393928
$ws.Cells.Item(9, 9).Value = 0.5569842765993244
$ws.Cells.Item(9, 10).Value = 0.5569842765993244
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.585130333333333
$ws.Cells.Item(9, 14).Value = 7.755391
$ws.Cells.Item(9, 15).Value = 0.1971906014874617
$ws.Cells.Item(9, 16).Value = 0.1971906014874618
$ws.Cells.Item(9, 17).Value = 40.839826962872
$ws.Cells.Item(9, 18).Value = 367.558442665848
$ws.Cells.Item(9, 19).Value = 0.1098320645216795
$ws.Cells.Item(9, 20).Value = 0.1098320645216795

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Tfpi"
$ws.Cells.Item(10, 3).Value = "Vldlr"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.587063
$ws.Cells.Item(10, 8).Value = 10.761189
$ws.Cells.Item(10, 9).Value = 0.1264679532473782
$ws.Cells.Item(10, 10).Value = 0.1264679532473782
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.463191
$ws.Cells.Item(10, 14).Value = 1.389573
$ws.Cells.Item(10, 15).Value = 0.0353316468093919
$ws.Cells.Item(10, 16).Value = 0.0353316468093919
$ws.Cells.Item(10, 17).Value = 1.661495298033
$ws.Cells.Item(10, 18).Value = 14.953457682297
$ws.Cells.Item(10, 19).Value = 0.004468321056843053
$ws.Cells.Item(10, 20).Value = 0.004468321056843053

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Tfpi"
$ws.Cells.Item(11, 3).Value = "Vldlr"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 3.587063
$ws.Cells.Item(11, 8).Value = 10.761189
$ws.Cells.Item(11, 9).Value = 0.1264679532473782
$ws.Cells.Item(11, 10).Value = 0.1264679532473782
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 9.932929333333332
$ws.Cells.Item(11, 14).Value = 29.798788
$ws.Cells.Item(11, 15).Value = 0.7576717833204485
$ws.Cells.Item(11, 16).Value = 0.7576717833204486
$ws.Cells.Item(11, 17).Value = 35.63004329321466
$ws.Cells.Item(11, 18).Value = 320.670389638932
$ws.Cells.Item(11, 19).Value = 0.09582119966982812
$ws.Cells.Item(11, 20).Value = 0.09582119966982813

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Tfpi"
$ws.Cells.Item(12, 3).Value = "Vldlr"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 3.587063
$ws.Cells.Item(12, 8).Value = 10.761189
$ws.Cells.Item(12, 9).Value = 0.1264679532473782
$ws.Cells.Item(12, 10).Value = 0.1264679532473782
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.1285543333333333
$ws.Cells.Item(12, 14).Value = 0.385663
$ws.Cells.Item(12, 15).Value = 0.009805968382697785
$ws.Cells.Item(12, 16).Value = 0.009805968382697785
$ws.Cells.Item(12, 17).Value = 0.4611324925896667
$ws.Cells.Item(12, 18).Value = 4.150192433307
$ws.Cells.Item(12, 19).Value = 0.001240140750968292
$ws.Cells.Item(12, 20).Value = 0.001240140750968292

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Tfpi"
$ws.Cells.Item(13, 3).Value = "Vldlr"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 3.587063
$ws.Cells.Item(13, 8).Value = 10.761189
$ws.Cells.Item(13, 9).Value = 0.1264679532473782
$ws.Cells.Item(13, 10).Value = 0.1264679532473782
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.585130333333333
$ws.Cells.Item(13, 14).Value = 7.755391
$ws.Cells.Item(13, 15).Value = 0.1971906014874617
$ws.Cells.Item(13, 16).Value = 0.1971906014874618
$ws.Cells.Item(13, 17).Value = 9.273025368877667
$ws.Cells.Item(13, 18).Value = 83.457228319899
$ws.Cells.Item(13, 19).Value = 0.02493829176973869
$ws.Cells.Item(13, 20).Value = 0.02493829176973869

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Tfpi"
$ws.Cells.Item(14, 3).Value = "Vldlr"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 3.526644000000001
$ws.Cells.Item(14, 8).Value = 10.579932
$ws.Cells.Item(14, 9).Value = 0.1243377795461487
$ws.Cells.Item(14, 10).Value = 0.1243377795461487
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.463191
$ws.Cells.Item(14, 14).Value = 1.389573
$ws.Cells.Item(14, 15).Value = 0.0353316468093919
$ws.Cells.Item(14, 16).Value = 0.0353316468093919
$ws.Cells.Item(14, 17).Value = 1.633509761004
$ws.Cells.Item(14, 18).Value = 14.701587849036
$ws.Cells.Item(14, 19).Value = 0.004393058511988558
$ws.Cells.Item(14, 20).Value = 0.004393058511988558

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Tfpi"
$ws.Cells.Item(15, 3).Value = "Vldlr"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 3.526644000000001
$ws.Cells.Item(15, 8).Value = 10.579932
$ws.Cells.Item(15, 9).Value = 0.1243377795461487
$ws.Cells.Item(15, 10).Value = 0.1243377795461487
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 9.932929333333332
$ws.Cells.Item(15, 14).Value = 29.798788
$ws.Cells.Item(15, 15).Value = 0.7576717833204485
$ws.Cells.Item(15, 16).Value = 0.7576717833204486
$ws.Cells.Item(15, 17).Value = 35.029905635824
$ws.Cells.Item(15, 18).Value = 315.269150722416
$ws.Cells.Item(15, 19).Value = 0.09420722716283526
$ws.Cells.Item(15, 20).Value = 0.09420722716283528

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Tfpi"
$ws.Cells.Item(16, 3).Value = "Vldlr"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 3.526644000000001
$ws.Cells.Item(16, 8).Value = 10.579932
$ws.Cells.Item(16, 9).Value = 0.1243377795461487
$ws.Cells.Item(16, 10).Value = 0.1243377795461487
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1285543333333333
$ws.Cells.Item(16, 14).Value = 0.385663
$ws.Cells.Item(16, 15).Value = 0.009805968382697785
$ws.Cells.Item(16, 16).Value = 0.009805968382697785
$ws.Cells.Item(16, 17).Value = 0.453365368324
$ws.Cells.Item(16, 18).Value = 4.080288314916
$ws.Cells.Item(16, 19).Value = 0.001219252335004381
$ws.Cells.Item(16, 20).Value = 0.001219252335004381

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Tfpi"
$ws.Cells.Item(17, 3).Value = "Vldlr"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 3.526644000000001
$ws.Cells.Item(17, 8).Value = 10.579932
$ws.Cells.Item(17, 9).Value = 0.1243377795461487
$ws.Cells.Item(17, 10).Value = 0.1243377795461487
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 2.585130333333333
$ws.Cells.Item(17, 14).Value = 7.755391
$ws.Cells.Item(17, 15).Value = 0.1971906014874617
$ws.Cells.Item(17, 16).Value = 0.1971906014874618
$ws.Cells.Item(17, 17).Value = 9.116834379268001
$ws.Cells.Item(17, 18).Value = 82.05150941341202
$ws.Cells.Item(17, 19).Value = 0.02451824153632048
$ws.Cells.Item(17, 20).Value = 0.02451824153632048
